# Apply weekly update: insert a new record row at row 107 (pushing the
# existing rows 107-211 down to 108-212) and populate the new row with the
# newest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 107; this shifts rows
# 107..211 down to 108..212 while preserving their values/formatting.
$ws.Rows.Item(107).Insert()

# Fill in the newly inserted row 107 with the new data point.
$ws.Range("A107").Value = 4
$ws.Range("B107").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value = "Los Lagos"
$ws.Range("D107").Value = 44740
$ws.Range("E107").Value = 10
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100108
$ws.Range("H107").Value = "Tropicales y subtropicales"
$ws.Range("I107").Value = 100108002
$ws.Range("J107").Value = "Mango"
$ws.Range("K107").Value = "Sin especificar"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 200
$ws.Range("N107").Value = 9000
$ws.Range("O107").Value = 10000
$ws.Range("P107").Value = 9500
$ws.Range("Q107").Value = "$/bandeja 4 kilos"
$ws.Range("R107").Value = "Brasil"
$ws.Range("S107").Value = 2375
$ws.Range("T107").Value = 4

# Match the date cell's number format/style used by the rest of column D.
$ws.Range("D107").NumberFormat = $ws.Range("D108").NumberFormat
